# edgar6_v5 update N2O emissions all figures
#
# Updates the plot-data values (and their derived "label (xx%)" text)
# on all three sheets of the workbook to the refreshed EDGAR v6 numbers.

$wb = $excel.ActiveWorkbook

$wsDirect    = $wb.Worksheets.Item("direct emissions")
$wsIndirect  = $wb.Worksheets.Item("indirect emissions")
$wsSubsector = $wb.Worksheets.Item("indirect emissions - subsectors")

# ---------------------------------------------------------------------
# Sheet "direct emissions": column C (fraction, shown as %) rounded
# ---------------------------------------------------------------------
$wsDirect.Range("C4").Value = 23     # Electricity & heat   23.4 -> 23
$wsDirect.Range("C5").Value = 10     # Energy systems        10.4 -> 10
$wsDirect.Range("C6").Value = 24     # Industry               23.9 -> 24
$wsDirect.Range("C7").Value = 15     # Transport              14.8 -> 15

# ---------------------------------------------------------------------
# Sheet "indirect emissions": column E (fraction, shown as %) rounded
# ---------------------------------------------------------------------
$wsIndirect.Range("E2").Value  = 12    # Energy systems / GHG    12.4 -> 12
$wsIndirect.Range("E3").Value  = 12    # Energy systems / CO2_indirect
$wsIndirect.Range("E4").Value  = 34    # Industry / GHG          33.9 -> 34
$wsIndirect.Range("E5").Value  = 34    # Industry / CO2_indirect
$wsIndirect.Range("E8").Value  = 15    # Transport / GHG         15.1 -> 15
$wsIndirect.Range("E9").Value  = 15    # Transport / CO2_indirect
$wsIndirect.Range("E10").Value = 16    # Buildings / GHG         16.5 -> 16
$wsIndirect.Range("E11").Value = 16    # Buildings / CO2_indirect

# ---------------------------------------------------------------------
# Sheet "indirect emissions - subsectors": column E (fraction) and the
# matching column F label text ("<subsector> (<fraction>%)")
# ---------------------------------------------------------------------
$wsSubsector.Range("E10").Value = 13
$wsSubsector.Range("F10").Value = "Other (industry) (13%)"

$wsSubsector.Range("E11").Value = 0.12
$wsSubsector.Range("F11").Value = "Biomass burning (CH4, N2O) (0.12%)"

$wsSubsector.Range("E12").Value = 0.72
$wsSubsector.Range("F12").Value = "Manure management (N2O, CH4) (0.72%)"

$wsSubsector.Range("E13").Value = 0.75
$wsSubsector.Range("F13").Value = "Synthetic fertilizer application (N2O) (0.75%)"

$wsSubsector.Range("E17").Value = 11
$wsSubsector.Range("F17").Value = "LULUCF CO2 (11%)"

$wsSubsector.Range("E19").Value = 0.43
$wsSubsector.Range("F19").Value = "Rail  (0.43%)"

$wsSubsector.Range("E20").Value = 0.67
$wsSubsector.Range("F20").Value = "Domestic Aviation (0.67%)"

$wsSubsector.Range("E24").Value = 10
$wsSubsector.Range("F24").Value = "Road (10%)"

$wsSubsector.Range("E25").Value = 0.073
$wsSubsector.Range("F25").Value = "Non-CO2 (all buildings) (0.073%)"

$wsSubsector.Range("E27").Value = 11
$wsSubsector.Range("F27").Value = "Residential (11%)"
